# Refresh the cryptos price/volume table (Price + Volume(1h) columns), and
# swap the RenderToken/NEARProtocol rows (rows 37-38 traded rank places).
# Values are written with a leading "'" (quote-prefix) so Excel stores them
# as text, matching the sheet's original inline-string cell types and
# preserving exact formatting (thousand-separator dots, trailing zeros,
# padded "%" strings) instead of being auto-coerced to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.343.55"
$ws.Range('E2').Value = "'  -0.76%  "
$ws.Range('D3').Value = "'2.279.36"
$ws.Range('E3').Value = "'  -0.63%  "
$ws.Range('E4').Value = "'  -0.29%  "
$ws.Range('D5').Value = "'112.19"
$ws.Range('E5').Value = "'  -2.17%  "
$ws.Range('D6').Value = "'264.34"
$ws.Range('E6').Value = "'  -1.63%  "
$ws.Range('D7').Value = "'0.631"
$ws.Range('E7').Value = "'  +1.20%  "
$ws.Range('E8').Value = "'  +0.23%  "
$ws.Range('E9').Value = "'  -2.49%  "
$ws.Range('D10').Value = "'46.94"
$ws.Range('E10').Value = "'  -2.58%  "
$ws.Range('D11').Value = "'0.0934"
$ws.Range('E11').Value = "'  -0.39%  "
$ws.Range('D12').Value = "'9.25"
$ws.Range('E12').Value = "'  +5.12%  "
$ws.Range('E13').Value = "'  +1.61%  "
$ws.Range('D14').Value = "'15.43"
$ws.Range('E14').Value = "'  -1.26%  "
$ws.Range('D15').Value = "'2.625.59"
$ws.Range('E15').Value = "'  -0.44%  "
$ws.Range('D17').Value = "'2.285.58"
$ws.Range('E17').Value = "'  -0.29%  "
$ws.Range('D18').Value = "'43.166.49"
$ws.Range('E18').Value = "'  -1.21%  "
$ws.Range('E19').Value = "'  -1.53%  "
$ws.Range('D20').Value = "'6.76"
$ws.Range('E20').Value = "'  +3.79%  "
$ws.Range('D21').Value = "'71.87"
$ws.Range('E21').Value = "'  -0.86%  "
$ws.Range('E22').Value = "'  -0.85%  "
$ws.Range('D23').Value = "'234.39"
$ws.Range('E23').Value = "'  +0.51%  "
$ws.Range('D24').Value = "'9.46"
$ws.Range('E24').Value = "'  -2.97%  "
$ws.Range('E25').Value = "'  +0.84%  "
$ws.Range('E26').Value = "'  +1.72%  "
$ws.Range('D27').Value = "'11.36"
$ws.Range('E27').Value = "'  -2.64%  "
$ws.Range('D28').Value = "'3.98"
$ws.Range('E28').Value = "'  +0.05%  "
$ws.Range('D29').Value = "'40.36"
$ws.Range('E29').Value = "'  -3.74%  "
$ws.Range('D30').Value = "'3.34"
$ws.Range('E30').Value = "'  -1.60%  "
$ws.Range('E31').Value = "'  -1.01%  "
$ws.Range('D32').Value = "'173.08"
$ws.Range('E32').Value = "'  -1.92%  "
$ws.Range('D33').Value = "'21.49"
$ws.Range('E33').Value = "'  -0.42%  "
$ws.Range('D34').Value = "'0.0902"
$ws.Range('E34').Value = "'  -3.43%  "
$ws.Range('D35').Value = "'5.67"
$ws.Range('E35').Value = "'  +1.53%  "
$ws.Range('E36').Value = "'  +0.83%  "
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = "'4.07"
$ws.Range('E37').Value = "'  +5.92%  "
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = "'4.64"
$ws.Range('E38').Value = "'  -2.30%  "
$ws.Range('D39').Value = "'0.0370"
$ws.Range('E39').Value = "'  +3.40%  "
$ws.Range('D40').Value = "'0.104"
$ws.Range('E40').Value = "'  -4.70%  "
$ws.Range('D41').Value = "'2.61"
$ws.Range('E41').Value = "'  +6.88%  "
$ws.Range('D42').Value = "'76.33"
$ws.Range('E42').Value = "'  +4.52%  "
$ws.Range('D43').Value = "'13.99"
$ws.Range('E43').Value = "'  +1.21%  "
$ws.Range('E44').Value = "'  -2.56%  "
$ws.Range('D45').Value = "'6.12"
$ws.Range('E45').Value = "'  -0.54%  "
$ws.Range('D46').Value = "'1.00"
$ws.Range('E46').Value = "'  -0.06%  "
$ws.Range('D47').Value = "'1.37"
$ws.Range('E47').Value = "'  -4.76%  "
$ws.Range('D48').Value = "'103.53"
$ws.Range('E48').Value = "'  +0.57%  "
$ws.Range('D49').Value = "'8.59"
$ws.Range('E49').Value = "'  -1.46%  "
$ws.Range('E50').Value = "'  +2.54%  "
$ws.Range('D51').Value = "'0.0994"
$ws.Range('E51').Value = "'  -0.83%  "
